$d = $word.ActiveDocument

# 1. Ativação date update
$d.Content.Find.Execute(
    "Ativação: 01/01/2012", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ativação: 01/01/2025", 2) | Out-Null

# 2. Objetivos (Portuguese) - full replacement
$d.Content.Find.Execute(
    "Proporcionar ao estudante conhecimentos práticos sobre os a engenharia bioquímica, os processos bioquímicos e a enzimologia.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Desenvolver nos discentes as competências e habilidades necessárias para aplicar conhecimentos científicos, tecnológicos e de engenharia na concepção, projeto, instalação, otimização, supervisão e avaliação crítica da operação de bioprocessos, com ênfase em proporcionar ao estudante conhecimentos práticos sobre: 1) engenharia bioquímica, 2) processos bioquímicos e 3) enzimologia.",
    2) | Out-Null

# 3. Objetivos (English) - the run is currently empty (italic), insert text.
$para7 = $d.Paragraphs(7)
$para7.Range.InsertAfter("Develop in students the competencies and skills necessary to apply scientific, technological, and engineering knowledge in the design, project, installation, optimization, supervision, and critical evaluation of bioprocess operations, with an emphasis on providing students with practical knowledge in: 1) biochemical engineering, 2) biochemical processes, and 3) enzymology.")

# 4. Programa resumido (Portuguese) - full replacement
$d.Content.Find.Execute(
    "Condução de processos bioquímicos desde as etapas de preparação e esterilização de meio até as etapas de recuperação e caracterização de produtos",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "O processo engloba a montagem, esterilização e operação do biorreator, incluindo a calibração das sondas de pH e oxigênio dissolvido, além da determinação experimental do coeficiente volumétrico de transferência de massa e do tempo de mistura em reator de bancada. Também envolve a imobilização de levedura em esferas de alginato de cálcio e a quantificação das células imobilizadas, junto com a realização de um processo descontínuo no biorreator de bancada, monitorando variáveis de processo e analisando amostras para avaliar a concentração de células, substrato e produto. Adicionalmente, inclui atividades de purificação de enzimas por cromatografia líquida, caracterização da massa molar de enzimas e planejamento experimental para produção de etanol de segunda geração por hidrólise de celulose, com avaliação da viabilidade do processo.",
    2) | Out-Null

# 5. Programa resumido (English) - full replacement
$d.Content.Find.Execute(
    "Proceedings of biochemical processes since the steps of medium preparation and sterilization until the steps of products recuperation and characterization.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The process encompasses the assembly, sterilization, and operation of the bioreactor, including the calibration of pH and dissolved oxygen probes, as well as the experimental determination of the volumetric mass transfer coefficient and the mixing time in a bench-scale reactor. It also involves the immobilization of yeast in calcium alginate beads and the quantification of the immobilized cells, along with the execution of a batch process in the bench-scale bioreactor, monitoring process variables and analyzing samples to assess cell, substrate, and product concentrations. Additionally, it includes enzyme purification activities through liquid chromatography, characterization of enzyme molar mass, and experimental planning for the production of second-generation ethanol by cellulose hydrolysis, evaluating the feasibility of the process.",
    2) | Out-Null

# 6. Programa (Portuguese) - paragraph with multiple runs/breaks becomes a single run
$para14 = $d.Paragraphs(14)
$r14 = $d.Range($para14.Range.Start, $para14.Range.End - 1)
$r14.Text = "1) Montagem, esterilização e operação do biorreator, com calibração das sondas de pH e oxigênio dissolvido.2) Determinação experimental do coeficiente volumétrico de transferência de massa (kLa) usando o método `"Gassing out`" em diferentes condições de agitação e aeração.3) Determinação experimental do tempo de mistura em reator de bancada pelo método de descoloração.4) Imobilização de levedura em esferas de alginato de cálcio, seguida da quantificação das células imobilizadas.5) Realização de processo descontínuo em biorreator de bancada, monitorando variáveis de processo e analisando amostras para avaliar concentração de células, substrato e produto.6) Projeto de purificação de enzimas por cromatografia líquida através de equipamento do tipo FPLC onde um extrato enzimático  é utilizado como amostra teste. São utilizadas técnicas de cromatografia de troca iônica e também de exclusão de tamanho. Após, os grupos, recebem uma amostra pura e caracterizam esta do ponto de vista bioquímico (pH e temperatura ideal, estabilidade térmica), definindo quais os parâmetros ideais para aplicação.7) Caracterização de enzimas quanto a massa molar: calibração de uma coluna cromatográfica com proteínas conhecidas e determinação da massa molar de enzima problema; Determinação da enzima problema via absorção em 280 nm e por atividade específica.8) Planejamento experimental com os resultados obtidos para definir a melhor condição de hidrólise de celulose para a produção de etanol de segunda geração, avaliando a viabilidade do processo."

# 7. Programa (English) - paragraph with multiple runs/breaks becomes a single run
$para15 = $d.Paragraphs(15)
$r15 = $d.Range($para15.Range.Start, $para15.Range.End - 1)
$r15.Text = "1)Assembly, sterilization, and operation of the bioreactor, including the calibration of pH and dissolved oxygen probes.2)Experimental determination of the volumetric mass transfer coefficient (kLa) using the `"Gassing out`" method under different agitation and aeration conditions.3)Experimental determination of the mixing time in a bench-scale reactor using the decolorization method.4)Immobilization of yeast in calcium alginate beads, followed by the quantification of immobilized cells.5)Execution of a batch process in a bench-scale bioreactor, monitoring process variables and analyzing samples to assess cell, substrate, and product concentrations.6)Enzyme purification project using liquid chromatography with an FPLC system where an enzyme extract is used as a test sample. Techniques include ion exchange chromatography and size exclusion chromatography. Subsequently, the groups receive a pure sample and characterize it biochemically (optimal pH and temperature, thermal stability), defining the ideal parameters for application.7)Characterization of enzymes regarding molar mass: calibration of a chromatographic column with known proteins and determination of the molar mass of the enzyme in question; determination of the enzyme in question via absorption at 280 nm and specific activity.8)Experimental planning with the obtained results to define the best conditions for cellulose hydrolysis for the production of second-generation ethanol, evaluating the feasibility of the process."

# 8. Avaliação - Método
$d.Content.Find.Execute(
    "Relatórios e seminários sobre os experimentos",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Os alunos serão avaliados por relatórios (R) e seminários (S) sobre os experimentos. A ponderação das notas será de 70% para a média aritmética das notas dos relatórios (R) e 30% para a média aritmética das notas dos seminários (S), ou seja: Média Final do período letivo normal (MF) = (0,7xR +0,3xS).",
    2) | Out-Null

# 9. Avaliação - Critério
$d.Content.Find.Execute(
    "Média aritmética entre os relatórios e seminários",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Serão aprovados os alunos que obtiverem média do período letivo normal igual ou maior que 5.",
    2) | Out-Null

# 10. Avaliação - Norma de recuperação
$d.Content.Find.Execute(
    "A recuperação será feita por meio de prova escrita (PR) e a média final (MF) será calculada pela equação: MF = (NF + PR)/2.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "A recuperação será feita por meio de uma prova escrita (PR) para alunos que tenham MF maior ou igual a 3,0 e menor do que 5,0 e pelo menos 70% de frequência. A nota de recuperação (NR) será a média simples entre a média final (MF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0.",
    2) | Out-Null

# 11. Bibliografia - merge 5 lines (with <w:br/>) into a single run of text
$para19 = $d.Paragraphs(19)
$r19 = $d.Range($para19.Range.Start, $para19.Range.End - 1)
$r19.Text = "1. BORZANI, W., SCHMIDELL, W., LIMA, U.A., AQUARONE, E. Biotecnologia Industrial - Fundamentos (Vol 1). São Paulo: Edgard Blucher Ltda, 2001. 2. SCHMIDELL, W., LIMA, U.A., AQUARONE, E., BORZANI, W. Biotecnologia Industrial - Engenharia Bioquímica (Vol 2), São Paulo: Edgard Blucher Ltda, 2001. 3. COPELAND R.A. Enzymes: a practical introduction to structure, mechanism and data analysis, New York: Academic Press, 2000 4. BON, E.S., FERRARA, M.A., CORVO, M.L. (Eds.) Enzimas em Biotecnologia - Produção, aplicação e mercado, Rio de Janeiro: Editora Interciêcnia, 2008. 5. Godfrey T. e West S (eds), Industrial Enzymology, Chapman-Hall, New York, 1996"

# 12. Requisitos - swap order of the two requirement lines
$para21 = $d.Paragraphs(21)
$r21 = $d.Range($para21.Range.Start, $para21.Range.End - 1)
$r21.Text = "LOT2017 -  Enzimologia  (Requisito fraco)" + [char]11 + "LOT2013 -  Engenharia Bioquímica I  (Requisito fraco)" + [char]11
